$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(100).Insert()

$ws.Cells.Item(100, 1).Value = 9
$ws.Cells.Item(100, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(100, 3).Value = "Metropolitana"
$ws.Cells.Item(100, 4).Value = 44582
$ws.Cells.Item(100, 5).Value = 13
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100101
$ws.Cells.Item(100, 8).Value = "Berries"
$ws.Cells.Item(100, 9).Value = 100101001
$ws.Cells.Item(100, 10).Value = "Arándano (blue)"
$ws.Cells.Item(100, 11).Value = "Sin especificar"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 830
$ws.Cells.Item(100, 14).Value = 3500
$ws.Cells.Item(100, 15).Value = 4000
$ws.Cells.Item(100, 16).Value = 3789
$ws.Cells.Item(100, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(100, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(100, 19).Value = 1894
$ws.Cells.Item(100, 20).Value = 2
